$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells K1, L1 - copy style from existing header cell (A1) to get bold/border/alignment
$ws.Range("K1").Value = "fxppo2_accuracy_qkeras"
$ws.Range("L1").Value = "orig-fxppo2-drop_qkeras"
$ws.Range("A1").Copy()
$ws.Range("K1:L1").PasteSpecial(-4122) # xlPasteFormats

$data = @{
    2  = @(0.4882790165809034, 0.1226415094339623)
    3  = @(0.2001143510577473, 0.137221269296741)
    4  = @(0.2264150943396226, 0.245854774156661)
    5  = @(0.3776443682104059, 0.1638078902229846)
    6  = @(0.254145225843339, 0.319325328759291)
    7  = @(0.4485420240137221, 0.0920526014865638)
    8  = @(0.2650085763293311, 0.1503716409376786)
    9  = @(0.3662092624356775, 0.2075471698113207)
    10 = @(0.2884505431675243, 0.09748427672955984)
    11 = @(0.213264722698685, 0.09348198970840482)
    12 = @(0.2492853058890795, 0.09405374499714123)
    13 = @(0.4082332761578045, 0.1580903373356204)
    14 = @(0.3847913093196112, 0.2447112635791881)
    15 = @(0.4405374499714123, 0.1249285305889079)
    16 = @(0.4622641509433962, 0.1503716409376786)
    17 = @(0.1497998856489423, 0.1092052601486563)
    18 = @(0.440251572327044, 0.2166952544311035)
    19 = @(0.1703830760434534, 0.3078902229845626)
    20 = @(0.2384219554030875, 0.1623785020011435)
    21 = @(0.5343053173241853, 0.01000571755288737)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 11).Value = $vals[0]
    $ws.Cells.Item($r, 12).Value = $vals[1]
}
